$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from A47 into new rows A48:A49 before setting values
$ws.Range("A47").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A49").PasteSpecial(-4122)

# Update B and C values for rows 2-47, and populate new rows 48-49 (A, B, C)
$data = @(
  @(2, "", "2.436526759445814", "4.73868487938369"),
  @(3, "", "3.07778131355755", "9.030399598100548"),
  @(4, "", "3.861623001282576", "13.65123610159381"),
  @(5, "", "5.144260277174242", "18.68636893467628"),
  @(6, "", "9.091850423698734", "22.9903699430008"),
  @(7, "", "10.52768881164416", "27.67055998720142"),
  @(8, "", "14.40878583818409", "32.25587837970533"),
  @(9, "", "18.16973274804635", "36.80757980118673"),
  @(10, "", "20.86967393234801", "41.41922070219037"),
  @(11, "", "21.98868764415236", "46.01868879703605"),
  @(12, "", "26.47033935895532", "50.20438883734235"),
  @(13, "", "28.12166889803166", "54.50282149219203"),
  @(14, "", "29.36431219975522", "58.75094118523278"),
  @(15, "", "32.59198483043343", "63.44251568013605"),
  @(16, "", "37.49818118510792", "68.2562384906346"),
  @(17, "", "38.5704179864037", "73.17758598011956"),
  @(18, "", "39.9518111210485", "77.62985714399545"),
  @(19, "", "42.076861034495", "83.88583665251386"),
  @(20, "", "45.50608037691438", "88.32491025662799"),
  @(21, "", "46.94632498548744", "92.84828262083994"),
  @(22, "", "50.57039674103132", "97.66255477013813"),
  @(23, "", "52.74798902292862", "102.294461109647"),
  @(24, "", "53.82551200476517", "106.8815747183173"),
  @(25, "", "54.8557200898988", "111.7607582365735"),
  @(26, "", "56.12649452687561", "116.4278048907276"),
  @(27, "", "57.68418335795553", "120.8479506220677"),
  @(28, "", "58.8910464855331", "125.3812356148811"),
  @(29, "", "59.55545538932339", "129.7360179214687"),
  @(30, "", "61.47758263611268", "134.5296403974218"),
  @(31, "", "63.05848452388574", "139.4791176579361"),
  @(32, "", "65.29076685512884", "143.8733334876593"),
  @(33, "", "66.71755706330386", "148.5734577686894"),
  @(34, "", "68.2848906453424", "153.221066593537"),
  @(35, "", "71.01773294007006", "157.5720652048809"),
  @(36, "", "72.73063753830674", "162.0787090265635"),
  @(37, "", "73.29587515367113", "166.4171714311537"),
  @(38, "", "75.11049528708644", "170.6716151184"),
  @(39, "", "77.0307028725347", "176.0904183787233"),
  @(40, "", "78.64739955066166", "180.4809625412018"),
  @(41, "", "80.24383491266964", "184.9087643332743"),
  @(42, "", "82.13551439585373", "189.5673878522085"),
  @(43, "", "82.4283487520233", "194.1010802888335"),
  @(44, "", "84.07630018341753", "198.6031439909085"),
  @(45, "", "88.63539008843144", "202.9477375643535"),
  @(46, "", "92.67940803563664", "207.6541678844868"),
  @(47, "", "94.3615718397746", "212.4519646700301"),
  @(48, "46", "95.59146200405738", "216.8299721989482"),
  @(49, "47", "96.5567941184588", "221.7766896826433")
)

foreach ($entry in $data) {
  $r = $entry[0]
  $aVal = $entry[1]
  $bVal = [double]$entry[2]
  $cVal = [double]$entry[3]
  if ($aVal -ne "") {
    $ws.Cells.Item($r, 1).Value = [double]$aVal
  }
  $ws.Cells.Item($r, 2).Value = $bVal
  $ws.Cells.Item($r, 3).Value = $cVal
}
